$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that currently sits at the end
#    of the last paragraph ("... Eksistere"). It will be re-added later at
#    the end of the new final paragraph.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Build a range collapsed to just before the very end of the document
#    (i.e. right before the final paragraph mark) and insert the four new
#    paragraphs as raw WordprocessingML so the resulting markup matches the
#    target formatting exactly (style, spacing, theme color, bookmarks).
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="100" w:name="_Hlk34726438"/>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
    <w:t>Postconditions:</w:t>
  </w:r>
  <w:bookmarkEnd w:id="100"/>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
  </w:pPr>
  <w:r>
    <w:t>En instans ib af Indtjeningsbidrag blev oprettet</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
  </w:pPr>
  <w:r>
    <w:t>ib.beløb blev sat til KKO – mfb</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
  </w:pPr>
  <w:r>
    <w:t>ib blev presenteret for h</w:t>
  </w:r>
  <w:bookmarkStart w:id="101" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="101"/>
</w:p>
</w:body>
</w:document>
'@

$null = $endRange.InsertXML($xmlSnippet)

Write-Host "Edit applied successfully"
